$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text (not numeric) storage for Price column cells that look like numbers
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply updated values
$ws.Range("D2").Value = '39.807.78'
$ws.Range("E2").Value = '  +1.15%  '
$ws.Range("D3").Value = '2.211.34'
$ws.Range("E3").Value = '  +0.87%  '
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").Value = '291.91'
$ws.Range("E5").Value = '  -1.10%  '
$ws.Range("D6").Value = '86.57'
$ws.Range("E6").Value = '  +5.92%  '
$ws.Range("D7").Value = '0.514'
$ws.Range("E7").Value = '  +0.88%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("D9").Value = '0.471'
$ws.Range("E9").Value = '  +0.72%  '
$ws.Range("D10").Value = '30.21'
$ws.Range("E10").Value = '  +3.94%  '
$ws.Range("D11").Value = '0.0783'
$ws.Range("E11").Value = '  +2.01%  '
$ws.Range("D12").Value = '47.51'
$ws.Range("E12").Value = '  +1.04%  '
$ws.Range("D13").Value = '0.108'
$ws.Range("E13").Value = '  +1.69%  '
$ws.Range("D14").Value = '6.32'
$ws.Range("E14").Value = '  +1.52%  '
$ws.Range("D15").Value = '2.554.44'
$ws.Range("E15").Value = '  +0.85%  '
$ws.Range("D16").Value = '13.97'
$ws.Range("E16").Value = '  +0.50%  '
$ws.Range("D17").Value = '2.213.26'
$ws.Range("E17").Value = '  +0.79%  '
$ws.Range("D18").Value = '0.728'
$ws.Range("E18").Value = '  +2.68%  '
$ws.Range("D19").Value = '39.740.95'
$ws.Range("E19").Value = '  +1.30%  '
$ws.Range("D20").Value = '0.0₃0879'
$ws.Range("E20").Value = '  +1.41%  '
$ws.Range("D21").Value = '11.25'
$ws.Range("E21").Value = '  +9.53%  '
$ws.Range("D22").Value = '5.78'
$ws.Range("E22").Value = '  +1.29%  '
$ws.Range("D23").Value = '65.55'
$ws.Range("E23").Value = '  +1.28%  '
$ws.Range("D24").Value = '235.27'
$ws.Range("E24").Value = '  +4.54%  '
$ws.Range("E25").Value = '  -0.11%  '
$ws.Range("D26").Value = '2.46'
$ws.Range("E26").Value = '  +2.87%  '
$ws.Range("D27").Value = '1.83'
$ws.Range("E27").Value = '  +2.06%  '
$ws.Range("D28").Value = '22.68'
$ws.Range("E28").Value = '  +0.70%  '
$ws.Range("E29").Value = '  +3.88%  '
$ws.Range("D30").Value = '9.25'
$ws.Range("E30").Value = '  +2.00%  '
$ws.Range("D31").Value = '32.64'
$ws.Range("E31").Value = '  +2.66%  '
$ws.Range("D32").Value = '151.68'
$ws.Range("E32").Value = '  +2.15%  '
$ws.Range("D33").Value = '1.00'
$ws.Range("E33").Value = '  -0.05%  '
$ws.Range("D34").Value = '4.92'
$ws.Range("E34").Value = '  +2.70%  '
$ws.Range("D35").Value = '0.0716'
$ws.Range("E35").Value = '  +3.35%  '
$ws.Range("D36").Value = '2.37'
$ws.Range("E36").Value = '  +1.98%  '
$ws.Range("E37").Value = '  +1.82%  '
$ws.Range("D38").Value = '2.78'
$ws.Range("E38").Value = '  +6.02%  '
$ws.Range("D39").Value = '15.87'
$ws.Range("E39").Value = '  +3.84%  '
$ws.Range("D40").Value = '0.0984'
$ws.Range("E40").Value = '  +3.14%  '
$ws.Range("D41").Value = '1.70'
$ws.Range("E41").Value = '  +3.31%  '
$ws.Range("B42").Value = 'Maker'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D42").Value = '2.059.65'
$ws.Range("E42").Value = '  +8.37%  '
$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D43").Value = '3.78'
$ws.Range("E43").Value = '  +4.77%  '
$ws.Range("B44").Value = 'ApeXProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D44").Value = '2.11'
$ws.Range("E44").Value = '  +1.31%  '
$ws.Range("D45").Value = '0.0267'
$ws.Range("E45").Value = '  +3.53%  '
$ws.Range("D46").Value = '9.91'
$ws.Range("E46").Value = '  +10.95%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").Value = '17.68'
$ws.Range("E47").Value = '  +10.40%  '
$ws.Range("D48").Value = '2.61'
$ws.Range("E48").Value = '  +0.89%  '
$ws.Range("D49").Value = '2.433.18'
$ws.Range("E49").Value = '  +1.25%  '
$ws.Range("D50").Value = '70.91'
$ws.Range("E50").Value = '  -1.12%  '
$ws.Range("D51").Value = '88.80'
$ws.Range("E51").Value = '  +2.11%  '
